$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range('D2')
$r.NumberFormat = '@'
$r.Value = '25.865.47'
$r.ClearFormats()

$r = $ws.Range('D3')
$r.NumberFormat = '@'
$r.Value = '1.736.39'
$r.ClearFormats()

$r = $ws.Range('E3')
$r.NumberFormat = '@'
$r.Value = '  -0.40%  '
$r.ClearFormats()

$r = $ws.Range('D4')
$r.NumberFormat = '@'
$r.Value = '0.9998'
$r.ClearFormats()

$r = $ws.Range('E4')
$r.NumberFormat = '@'
$r.Value = '  +0.03%  '
$r.ClearFormats()

$r = $ws.Range('D5')
$r.NumberFormat = '@'
$r.Value = '242.61'
$r.ClearFormats()

$r = $ws.Range('E5')
$r.NumberFormat = '@'
$r.Value = '  +5.12%  '
$r.ClearFormats()

$r = $ws.Range('D6')
$r.NumberFormat = '@'
$r.Value = '1.0000'
$r.ClearFormats()

$r = $ws.Range('E6')
$r.NumberFormat = '@'
$r.Value = '  +0.02%  '
$r.ClearFormats()

$r = $ws.Range('D7')
$r.NumberFormat = '@'
$r.Value = '0.5201'
$r.ClearFormats()

$r = $ws.Range('E7')
$r.NumberFormat = '@'
$r.Value = '  -0.94%  '
$r.ClearFormats()

$r = $ws.Range('D8')
$r.NumberFormat = '@'
$r.Value = '0.2743'
$r.ClearFormats()

$r = $ws.Range('E9')
$r.NumberFormat = '@'
$r.Value = '  +0.04%  '
$r.ClearFormats()

$r = $ws.Range('D10')
$r.NumberFormat = '@'
$r.Value = '1.739.35'
$r.ClearFormats()

$r = $ws.Range('E10')
$r.NumberFormat = '@'
$r.Value = '  +0.25%  '
$r.ClearFormats()

$r = $ws.Range('D11')
$r.NumberFormat = '@'
$r.Value = '0.07186'
$r.ClearFormats()

$r = $ws.Range('E11')
$r.NumberFormat = '@'
$r.Value = '  +1.28%  '
$r.ClearFormats()

$r = $ws.Range('D12')
$r.NumberFormat = '@'
$r.Value = '14.99'
$r.ClearFormats()

$r = $ws.Range('E12')
$r.NumberFormat = '@'
$r.Value = '  -2.05%  '
$r.ClearFormats()

$r = $ws.Range('E13')
$r.NumberFormat = '@'
$r.Value = '  -0.10%  '
$r.ClearFormats()

$r = $ws.Range('E14')
$r.NumberFormat = '@'
$r.Value = '  +1.72%  '
$r.ClearFormats()

$r = $ws.Range('D15')
$r.NumberFormat = '@'
$r.Value = '77.31'
$r.ClearFormats()

$r = $ws.Range('E15')
$r.NumberFormat = '@'
$r.Value = '  -0.31%  '
$r.ClearFormats()

$r = $ws.Range('D16')
$r.NumberFormat = '@'
$r.Value = '0.9999'
$r.ClearFormats()

$r = $ws.Range('E16')
$r.NumberFormat = '@'
$r.Value = '  +0.02%  '
$r.ClearFormats()

$r = $ws.Range('D17')
$r.NumberFormat = '@'
$r.Value = '0.9997'
$r.ClearFormats()

$r = $ws.Range('E17')
$r.NumberFormat = '@'
$r.Value = '  +0.03%  '
$r.ClearFormats()

$r = $ws.Range('D18')
$r.NumberFormat = '@'
$r.Value = '25.900.90'
$r.ClearFormats()

$r = $ws.Range('E18')
$r.NumberFormat = '@'
$r.Value = '  +0.02%  '
$r.ClearFormats()

$r = $ws.Range('D19')
$r.NumberFormat = '@'
$r.Value = '11.77'
$r.ClearFormats()

$r = $ws.Range('D20')
$r.NumberFormat = '@'
$r.Value = '0.000006778'
$r.ClearFormats()

$r = $ws.Range('E20')
$r.NumberFormat = '@'
$r.Value = '  +1.43%  '
$r.ClearFormats()

$r = $ws.Range('D21')
$r.NumberFormat = '@'
$r.Value = '1.962.24'
$r.ClearFormats()

$r = $ws.Range('E21')
$r.NumberFormat = '@'
$r.Value = '  -0.09%  '
$r.ClearFormats()

$r = $ws.Range('D22')
$r.NumberFormat = '@'
$r.Value = '4.281'
$r.ClearFormats()

$r = $ws.Range('E22')
$r.NumberFormat = '@'
$r.Value = '  -0.51%  '
$r.ClearFormats()

$r = $ws.Range('D23')
$r.NumberFormat = '@'
$r.Value = '8.629'
$r.ClearFormats()

$r = $ws.Range('E23')
$r.NumberFormat = '@'
$r.Value = '  -1.71%  '
$r.ClearFormats()

$r = $ws.Range('D24')
$r.NumberFormat = '@'
$r.Value = '5.289'
$r.ClearFormats()

$r = $ws.Range('E24')
$r.NumberFormat = '@'
$r.Value = '  +2.46%  '
$r.ClearFormats()

$r = $ws.Range('E25')
$r.NumberFormat = '@'
$r.Value = '  -2.29%  '
$r.ClearFormats()

$r = $ws.Range('D26')
$r.NumberFormat = '@'
$r.Value = '1.511'
$r.ClearFormats()

$r = $ws.Range('E26')
$r.NumberFormat = '@'
$r.Value = '  -0.56%  '
$r.ClearFormats()

$r = $ws.Range('D27')
$r.NumberFormat = '@'
$r.Value = '15.21'
$r.ClearFormats()

$r = $ws.Range('E27')
$r.NumberFormat = '@'
$r.Value = '  +0.26%  '
$r.ClearFormats()

$r = $ws.Range('D28')
$r.NumberFormat = '@'
$r.Value = '1.774'
$r.ClearFormats()

$r = $ws.Range('E28')
$r.NumberFormat = '@'
$r.Value = '  -1.29%  '
$r.ClearFormats()

$r = $ws.Range('D29')
$r.NumberFormat = '@'
$r.Value = '105.12'
$r.ClearFormats()

$r = $ws.Range('E29')
$r.NumberFormat = '@'
$r.Value = '  +2.43%  '
$r.ClearFormats()

$r = $ws.Range('D30')
$r.NumberFormat = '@'
$r.Value = '3.953'
$r.ClearFormats()

$r = $ws.Range('E30')
$r.NumberFormat = '@'
$r.Value = '  +6.02%  '
$r.ClearFormats()

$r = $ws.Range('D31')
$r.NumberFormat = '@'
$r.Value = '0.08254'
$r.ClearFormats()

$r = $ws.Range('E31')
$r.NumberFormat = '@'
$r.Value = '  -0.73%  '
$r.ClearFormats()

$r = $ws.Range('D32')
$r.NumberFormat = '@'
$r.Value = '3.651'
$r.ClearFormats()

$r = $ws.Range('E32')
$r.NumberFormat = '@'
$r.Value = '  +3.06%  '
$r.ClearFormats()

$r = $ws.Range('D33')
$r.NumberFormat = '@'
$r.Value = '0.04668'
$r.ClearFormats()

$r = $ws.Range('E33')
$r.NumberFormat = '@'
$r.Value = '  +3.09%  '
$r.ClearFormats()

$r = $ws.Range('E34')
$r.NumberFormat = '@'
$r.Value = '  +1.77%  '
$r.ClearFormats()

$r = $ws.Range('E35')
$r.NumberFormat = '@'
$r.Value = '  +1.22%  '
$r.ClearFormats()

$r = $ws.Range('D36')
$r.NumberFormat = '@'
$r.Value = '0.6201'
$r.ClearFormats()

$r = $ws.Range('E36')
$r.NumberFormat = '@'
$r.Value = '  -0.10%  '
$r.ClearFormats()

$r = $ws.Range('D37')
$r.NumberFormat = '@'
$r.Value = '2.691'
$r.ClearFormats()

$r = $ws.Range('E37')
$r.NumberFormat = '@'
$r.Value = '  +0.31%  '
$r.ClearFormats()

$r = $ws.Range('D38')
$r.NumberFormat = '@'
$r.Value = '0.01601'
$r.ClearFormats()

$r = $ws.Range('E38')
$r.NumberFormat = '@'
$r.Value = '  +0.70%  '
$r.ClearFormats()

$r = $ws.Range('D39')
$r.NumberFormat = '@'
$r.Value = '1.922'
$r.ClearFormats()

$r = $ws.Range('E39')
$r.NumberFormat = '@'
$r.Value = '  -0.91%  '
$r.ClearFormats()

$r = $ws.Range('E40')
$r.NumberFormat = '@'
$r.Value = '  -0.01%  '
$r.ClearFormats()

$r = $ws.Range('D41')
$r.NumberFormat = '@'
$r.Value = '99.82'
$r.ClearFormats()

$r = $ws.Range('E41')
$r.NumberFormat = '@'
$r.Value = '  -0.36%  '
$r.ClearFormats()

$r = $ws.Range('D42')
$r.NumberFormat = '@'
$r.Value = '0.3857'
$r.ClearFormats()

$r = $ws.Range('E42')
$r.NumberFormat = '@'
$r.Value = '  -0.46%  '
$r.ClearFormats()

$r = $ws.Range('D43')
$r.NumberFormat = '@'
$r.Value = '0.7464'
$r.ClearFormats()

$r = $ws.Range('E43')
$r.NumberFormat = '@'
$r.Value = '  +1.79%  '
$r.ClearFormats()

$r = $ws.Range('D44')
$r.NumberFormat = '@'
$r.Value = '5.001'
$r.ClearFormats()

$r = $ws.Range('E44')
$r.NumberFormat = '@'
$r.Value = '  -0.18%  '
$r.ClearFormats()

$r = $ws.Range('D45')
$r.NumberFormat = '@'
$r.Value = '0.1126'
$r.ClearFormats()

$r = $ws.Range('E45')
$r.NumberFormat = '@'
$r.Value = '  +0.11%  '
$r.ClearFormats()

$r = $ws.Range('D46')
$r.NumberFormat = '@'
$r.Value = '6.260'
$r.ClearFormats()

$r = $ws.Range('E46')
$r.NumberFormat = '@'
$r.Value = '  -0.22%  '
$r.ClearFormats()

$r = $ws.Range('D47')
$r.NumberFormat = '@'
$r.Value = '55.12'
$r.ClearFormats()

$r = $ws.Range('E47')
$r.NumberFormat = '@'
$r.Value = '  +3.01%  '
$r.ClearFormats()

$r = $ws.Range('D48')
$r.NumberFormat = '@'
$r.Value = '0.05210'
$r.ClearFormats()

$r = $ws.Range('E48')
$r.NumberFormat = '@'
$r.Value = '  -2.48%  '
$r.ClearFormats()

$r = $ws.Range('E49')
$r.NumberFormat = '@'
$r.Value = '  +1.51%  '
$r.ClearFormats()

$r = $ws.Range('D50')
$r.NumberFormat = '@'
$r.Value = '7.528'
$r.ClearFormats()

$r = $ws.Range('E50')
$r.NumberFormat = '@'
$r.Value = '  -1.95%  '
$r.ClearFormats()

$r = $ws.Range('D51')
$r.NumberFormat = '@'
$r.Value = '0.3413'
$r.ClearFormats()

$r = $ws.Range('E51')
$r.NumberFormat = '@'
$r.Value = '  -0.36%  '
$r.ClearFormats()

